# Elimina los EC (estados de cuenta) anteriores y agrega los nuevos;
# se actualiza la base de datos de trabajadores / periodos en mora.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room: insert 3 new data rows right after the current last data
# row (old row 31) so the table grows from 16 rows (16-31) to 19 rows
# (16-34); this also pushes the footer (signature) rows down from 36/37
# to 39/40. Copy number formats from the row above so the new rows match
# the existing "middle" row styling, leaving the final row's special
# (bottom-border) styling on row 34.
$ws.Rows(31).EntireRow.Insert()
$ws.Rows(31).EntireRow.Insert()
$ws.Rows(31).EntireRow.Insert()

$ws.Range("B30:J30").Copy()
$ws.Range("B31:J33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Replace the whole worker/periods table (rows 16-34) with the new
# dataset.
$ws.Range("B16").Value2 = "CC"
$ws.Range("C16").Value2 = "28215702"
$ws.Range("D16").Value2 = "NANCY SANDOVAL SALINAS"
$ws.Range("E16").Value2 = "2507"
$ws.Range("F16").Value2 = 2899
$ws.Range("G16").Value2 = 2173914

$ws.Range("B17").Value2 = "CC"
$ws.Range("C17").Value2 = "28215702"
$ws.Range("D17").Value2 = "NANCY SANDOVAL SALINAS"
$ws.Range("E17").Value2 = "2505"
$ws.Range("F17").Value2 = 2899
$ws.Range("G17").Value2 = 2173914

$ws.Range("B18").Value2 = "CC"
$ws.Range("C18").Value2 = "1104872748"
$ws.Range("D18").Value2 = "PEDRO LUIS PADILLA ARRIETA"
$ws.Range("E18").Value2 = "2303"
$ws.Range("F18").Value2 = 46400
$ws.Range("G18").Value2 = 1400000

$ws.Range("B19").Value2 = "CC"
$ws.Range("C19").Value2 = "1100394158"
$ws.Range("D19").Value2 = "IVAN JAVIER DOMINGUEZ HERNANDEZ"
$ws.Range("E19").Value2 = "2303"
$ws.Range("F19").Value2 = 56000
$ws.Range("G19").Value2 = 1500000

$ws.Range("B20").Value2 = "CC"
$ws.Range("C20").Value2 = "1100400317"
$ws.Range("D20").Value2 = "EDUAR ANTONIO CORREA CASTRO"
$ws.Range("E20").Value2 = "2304"
$ws.Range("F20").Value2 = 56000
$ws.Range("G20").Value2 = 1400000

$ws.Range("B21").Value2 = "CC"
$ws.Range("C21").Value2 = "1100548932"
$ws.Range("D21").Value2 = "RAFAEL EDUARDO VANEGAS LUNA"
$ws.Range("E21").Value2 = "2303"
$ws.Range("F21").Value2 = 56000
$ws.Range("G21").Value2 = 1400000

$ws.Range("B22").Value2 = "CC"
$ws.Range("C22").Value2 = "1100398349"
$ws.Range("D22").Value2 = "EMBER GUSTAVO DE LA OSSA ATENCIA"
$ws.Range("E22").Value2 = "2303"
$ws.Range("F22").Value2 = 46400
$ws.Range("G22").Value2 = 1300000

$ws.Range("B23").Value2 = "CC"
$ws.Range("C23").Value2 = "1100396230"
$ws.Range("D23").Value2 = "JESUS ALFONSO PINEDA AVILA"
$ws.Range("E23").Value2 = "2303"
$ws.Range("F23").Value2 = 46400
$ws.Range("G23").Value2 = 1160000

$ws.Range("B24").Value2 = "CC"
$ws.Range("C24").Value2 = "73242423"
$ws.Range("D24").Value2 = "LUIS CARLOS MENDEZ JIMENEZ"
$ws.Range("E24").Value2 = "2303"
$ws.Range("F24").Value2 = 1547
$ws.Range("G24").Value2 = 1160000

$ws.Range("B25").Value2 = "CC"
$ws.Range("C25").Value2 = "1052954826"
$ws.Range("D25").Value2 = "ALEXANDRA ISABEL VALLE BARRAGAN"
$ws.Range("E25").Value2 = "2303"
$ws.Range("F25").Value2 = 26000
$ws.Range("G25").Value2 = 908526

$ws.Range("B26").Value2 = "CC"
$ws.Range("C26").Value2 = "1100550563"
$ws.Range("D26").Value2 = "WINDER JOSE SOLORZANO RICARDO"
$ws.Range("E26").Value2 = "2303"
$ws.Range("F26").Value2 = 1547
$ws.Range("G26").Value2 = 1160000

$ws.Range("B27").Value2 = "CC"
$ws.Range("C27").Value2 = "9143691"
$ws.Range("D27").Value2 = "EDUARDO ALVARADO QUESADA"
$ws.Range("E27").Value2 = "2303"
$ws.Range("F27").Value2 = 60000
$ws.Range("G27").Value2 = 1500000

$ws.Range("B28").Value2 = "CC"
$ws.Range("C28").Value2 = "72333328"
$ws.Range("D28").Value2 = "GEIGER ENRRIQUE HERRERA CASTRO"
$ws.Range("E28").Value2 = "2303"
$ws.Range("F28").Value2 = 60000
$ws.Range("G28").Value2 = 1500000

$ws.Range("B29").Value2 = "CC"
$ws.Range("C29").Value2 = "72333328"
$ws.Range("D29").Value2 = "GEIGER ENRRIQUE HERRERA CASTRO"
$ws.Range("E29").Value2 = "2302"
$ws.Range("F29").Value2 = 60000
$ws.Range("G29").Value2 = 1500000

$ws.Range("B30").Value2 = "CC"
$ws.Range("C30").Value2 = "1005418669"
$ws.Range("D30").Value2 = "EDGAR DANILO PALENCIA ATENCIA"
$ws.Range("E30").Value2 = "2303"
$ws.Range("F30").Value2 = 46400
$ws.Range("G30").Value2 = 1300000

$ws.Range("B31").Value2 = "CC"
$ws.Range("C31").Value2 = "23221622"
$ws.Range("D31").Value2 = "ODALIS DEL CARMEN SOLIS HUERTAS"
$ws.Range("E31").Value2 = "2507"
$ws.Range("F31").Value2 = 56940
$ws.Range("G31").Value2 = 1423500

$ws.Range("B32").Value2 = "CC"
$ws.Range("C32").Value2 = "23221622"
$ws.Range("D32").Value2 = "ODALIS DEL CARMEN SOLIS HUERTAS"
$ws.Range("E32").Value2 = "2506"
$ws.Range("F32").Value2 = 56940
$ws.Range("G32").Value2 = 1423500

$ws.Range("B33").Value2 = "CC"
$ws.Range("C33").Value2 = "23221622"
$ws.Range("D33").Value2 = "ODALIS DEL CARMEN SOLIS HUERTAS"
$ws.Range("E33").Value2 = "2505"
$ws.Range("F33").Value2 = 56940
$ws.Range("G33").Value2 = 1423500

$ws.Range("B34").Value2 = "CC"
$ws.Range("C34").Value2 = "23221622"
$ws.Range("D34").Value2 = "ODALIS DEL CARMEN SOLIS HUERTAS"
$ws.Range("E34").Value2 = "2504"
$ws.Range("F34").Value2 = 28470
$ws.Range("G34").Value2 = 1423500

# --- Update the summary header figures above the table.
$ws.Range("E11").Value2 = 767782   # VALOR MORA total
$ws.Range("C13").Value2 = 14       # Cant. Trabajadores
$ws.Range("F13").Value2 = 7        # Cant. Periodos
